$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.533.37"
$ws.Range("E2").Value = "  -2.69%  "

$ws.Range("D3").Value = "2.616.95"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.37"
$ws.Range("E5").Value = "  -1.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.42"
$ws.Range("E6").Value = "  -2.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D9").Value = "2.641.91"
$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.35"
$ws.Range("E10").Value = "  +0.44%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  -1.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("E13").Value = "  -1.79%  "

$ws.Range("D14").Value = "3.093.62"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "58.704.64"
$ws.Range("E15").Value = "  -2.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.93"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "2.652.73"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  -2.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.78"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("E21").Value = "  -0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.07"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.06"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "2.755.84"
$ws.Range("E26").Value = "  +0.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.161"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  -0.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.41"
$ws.Range("E32").Value = "  +7.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.81"
$ws.Range("E33").Value = "  -0.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "150.02"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  +11.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.00"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.15"
$ws.Range("E38").Value = "  +1.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.849"
$ws.Range("E39").Value = "  -1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.45"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.67"
$ws.Range("E41").Value = "  +1.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.40"
$ws.Range("E42").Value = "  -1.41%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.614"
$ws.Range("E43").Value = "  -2.31%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "278.17"
$ws.Range("E44").Value = "  -4.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  -0.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0982"
$ws.Range("E46").Value = "  -1.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.48"
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0535"
$ws.Range("E48").Value = "  -1.93%  "

$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -1.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.66"
$ws.Range("E51").Value = "  +0.12%  "
